# Update the station-load figures for hours 1-17 (rows 2-18) and drop the
# trailing hours 18-24 (rows 19-25), which the source data no longer covers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1,  17218, 4744, 7729, 0, 4871, 0, 14088, 17074),
    @(2,  16791, 4524, 7518, 0, 4673, 0, 13227, 16352),
    @(3,  16178, 4318, 7339, 0, 4523, 0, 13023, 15877),
    @(4,  15705, 4269, 7268, 0, 4450, 0, 12784, 15796),
    @(5,  16704, 4560, 7432, 0, 4720, 0, 13319, 16447),
    @(6,  17753, 4910, 8062, 0, 5584, 0, 14034, 17459),
    @(7,  18187, 4410, 8441, 0, 6771, 0, 14185, 17577),
    @(8,  21165, 4679, 10305, 0, 9453, 0, 16193, 19363),
    @(9,  26539, 5435, 12881, 0, 12768, 0, 18743, 22087),
    @(10, 28816, 5831, 16459, 0, 13796, 0, 19437, 24090),
    @(11, 29162, 5801, 18360, 0, 14092, 0, 19402, 24922),
    @(12, 29265, 5919, 18421, 0, 14471, 0, 19716, 25305),
    @(13, 29689, 5765, 18171, 0, 14225, 0, 19981, 24921),
    @(14, 30436, 6044, 18805, 0, 15175, 0, 20615, 25900),
    @(15, 30921, 6017, 18794, 0, 14937, 0, 20903, 26070),
    @(16, 30681, 6080, 18729, 0, 14970, 0, 0,     25962),
    @(17, 29701, 6044, 18564, 0, 15022, 0, 0,     25317)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 1; $c -le 9; $c++) {
        $ws.Cells.Item($row, $c).Value = $vals[$c - 1]
    }
}

# The old rows 19-25 (hours 18-24) are gone entirely in the new data set.
$ws.Range("A19:I25").EntireRow.Delete()
